$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3:22 down to 4:23
$ws.Rows.Item(3).Insert()

# Copy formatting (border/font/alignment) from the row below onto the new row's label cell
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Populate the new row 3 with its label and data
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"
$ws.Range("B3").Value = 2.174397541324862
$ws.Range("C3").Value = 10.04419022232486
$ws.Range("D3").Value = -8.284815777675137
$ws.Range("E3").Value = -0.1447337776751375
$ws.Range("F3").Value = 1.643575222324863
$ws.Range("G3").Value = -1.534543777675137
$ws.Range("H3").Value = -1.829872777675137
